$d = $word.ActiveDocument

# --- First paragraph: add a paragraph border (space-only, no line), ---
# --- change left indent from 120 -> 225 twips, update the ID text,  ---
# --- and drop the trailing " " run entirely.                        ---

$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat

# Paragraph border: w:pBdr with top/left/bottom/right each carrying
# only w:space="5" (no line style/width/color).
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromRight = 5

# Left indent: 120 -> 225 twips (LeftIndent is expressed in points,
# 1 point = 20 twips, so 225/20 = 11.25).
$pf.LeftIndent = 11.25

# Remove the trailing " " run (second run in the paragraph) before the
# paragraph mark. The paragraph's range is
# "**ID__AFFARS_5301_topic_6__ID** " + paragraph mark.
$idLen = "**ID__AFFARS_5301_topic_6__ID**".Length
$spaceRange = $d.Range($p1.Range.Start + $idLen, $p1.Range.Start + $idLen + 1)
$spaceRange.Delete()

# Update the ID placeholder text on the (now sole) run in the paragraph.
$d.Content.Find.Execute("**ID__AFFARS_5301_topic_6__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5301_2__ID**", 2)
